# table_general_sums.xlsx -- "ajuste final nos graficos do apendice"
#
# Adds a second worksheet ("Planilha1") that repeats the appendix table with
# the servants/appointees/toplevel column order reversed (and a merged
# two-group header above the existing column headers), and tidies up the
# selection/tab state on the first sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- sheet1: selection / active state tidy-up -----------------------------
[void]$ws1.Range("A1:M13").Select()
$ws1.Rows(1).RowHeight = 42.75

# --- add the new sheet right after sheet1 ----------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Planilha1"

# --- row 2: column headers (same 12 metrics, servants/appointees/toplevel
#            order reversed within each of the 4 groups) ------------------
$ws2.Range("A2:M2").WrapText = $true
$ws2.Rows(2).RowHeight = 42.75

$ws2.Cells.Item(2, 1).Value = "year"
$ws2.Cells.Item(2, 2).Value = "servants 3 Ras"
$ws2.Cells.Item(2, 3).Value = "appointees_3RAs"
$ws2.Cells.Item(2, 4).Value = "toplevel_3RAs"
$ws2.Cells.Item(2, 5).Value = "servants_otherRAs"
$ws2.Cells.Item(2, 6).Value = "appointees_otherRAs"
$ws2.Cells.Item(2, 7).Value = "toplevel_otherRAs"
$ws2.Cells.Item(2, 8).Value = "servants_officers_3RAs"
$ws2.Cells.Item(2, 9).Value = "appointees_officers_3RAs"
$ws2.Cells.Item(2, 10).Value = "toplevel_officers_3RAs"
$ws2.Cells.Item(2, 11).Value = "servants_officers_otherRAs"
$ws2.Cells.Item(2, 12).Value = "appointees_officers_otherRAs"
$ws2.Cells.Item(2, 13).Value = "toplevel_officers_otherRAs"

# --- row 1: merged group headers over B:G and H:M -------------------------
$ws2.Cells.Item(1, 2).Value = "total (considering civil and military)"
$ws2.Cells.Item(1, 8).Value = "total (considering only military)"

$ws2.Range("B1:G1").HorizontalAlignment = -4108
$ws2.Range("H1:M1").HorizontalAlignment = -4108
$ws2.Range("H1:M1").WrapText = $true

[void]$ws2.Range("B1:G1").Merge()
[void]$ws2.Range("H1:M1").Merge()

# --- data rows 3-14: years 2013-2024, columns reordered vs. sheet1 --------
$years = 2013..2024
$data = @(
  ,@(4029,1106,332,5146,1640,601,28,25,10,9,5,5)
  ,@(4191,1192,356,5309,1776,650,34,20,10,16,4,4)
  ,@(4406,1222,370,5505,1794,671,42,24,10,22,4,4)
  ,@(4336,1245,380,5633,1903,701,39,23,9,29,6,4)
  ,@(4309,1327,398,5545,1932,738,36,22,8,27,6,3)
  ,@(4271,1307,367,5526,1939,747,41,17,7,30,6,4)
  ,@(4129,1327,373,5343,1921,759,35,13,6,27,8,5)
  ,@(3893,1287,377,6197,2158,800,40,21,13,26,9,6)
  ,@(4119,1336,372,6595,2110,812,53,30,20,37,10,8)
  ,@(4120,1320,376,6401,2130,827,51,31,22,28,9,7)
  ,@(4112,1336,392,6256,2217,892,46,31,22,32,12,8)
  ,@(4084,1343,400,6250,2146,875,17,16,14,4,4,2)
)

for ($i = 0; $i -lt $years.Count; $i++) {
  $r = $i + 3
  $ws2.Cells.Item($r, 1).Value = $years[$i]
  $ws2.Range($ws2.Cells.Item($r, 1), $ws2.Cells.Item($r, 1)).HorizontalAlignment = -4131
  $row = $data[$i]
  for ($c = 0; $c -lt $row.Count; $c++) {
    $ws2.Cells.Item($r, $c + 2).Value = $row[$c]
  }
}

[void]$ws2.Range("F11").Select()

Write-Output "done"
